# Iterationsplan.xlsx edit script
# Commit: "Lade till de olika symbolerna som används som pjäser"
#
# Summary of the change:
#  - Table 1 (rows 4-13) "iteration 0" plan: the first data row's iteration
#    label cell (A5) switches from the literal text "-1, Uppstart" to the
#    plain number -1 (the shared string "-1, Uppstart" becomes unused and is
#    dropped on save).
#  - Table 2 (rows 17-24), the status table, gets a brand-new leading column
#    "Iterationsnamn" (holding the iteration number, 0) and a new "Krav"
#    (requirement) column, which pushes the previous Uppgift/Status/
#    Uppskattad tid/Verklig tid columns one step to the right (B->C, C->D,
#    D->E, E->F). The "Formatering" row's status also flips from
#    "Ej påbörjat" to "Färdig".
#  - A few formatting-only marker cells (A25/A27/A29/A31) are added below
#    the table, the sheet's used range grows to F31, the selected cell
#    moves to D1, and column C is widened (it now holds long requirement
#    text instead of short status words).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlLeft = -4131

function Set-CellStyle($rng, [bool]$bold, [bool]$left) {
    if ($bold) { $rng.Font.Bold = $true }
    if ($left) { $rng.HorizontalAlignment = $xlLeft }
}

# ---------------------------------------------------------------------
# Table 1 (rows 4-13): A5 becomes a real number instead of quoted text.
# ---------------------------------------------------------------------
$ws.Range("A5").Value2 = -1

# ---------------------------------------------------------------------
# Table 2 (rows 17-24): insert "Iterationsnamn" + "Krav" columns, shifting
# the previous B:E content to C:F, and update the "Formatering" status.
# ---------------------------------------------------------------------

# -- Header row 17 --
$ws.Range("F17").Value2 = $ws.Range("E17").Value2   # Verklig tid
$ws.Range("E17").Value2 = $ws.Range("D17").Value2   # Uppskattad tid
$ws.Range("D17").Value2 = $ws.Range("C17").Value2   # Status
$ws.Range("C17").Value2 = $ws.Range("B17").Value2   # Uppgift
$ws.Range("B17").Value2 = $ws.Range("A17").Value2   # Krav (was "Krav" header already in A17)
$ws.Range("A17").Value2 = "Iterationsnamn"
Set-CellStyle $ws.Range("A17") $true $true

# -- Row 18 (iteration 0 / Planering) --
$ws.Range("F18").Value2 = $ws.Range("E18").Value2   # 3 (Verklig tid)
$ws.Range("E18").Value2 = $ws.Range("D18").Value2   # 3 (Uppskattad tid)
$ws.Range("D18").Value2 = $ws.Range("C18").Value2   # Påbörjat (Status)
$ws.Range("C18").Value2 = $ws.Range("B18").Value2   # Planering (Uppgift)
$ws.Range("B18").ClearContents() | Out-Null         # Krav column left blank
$ws.Range("A18").Value2 = 0                         # iteration number
Set-CellStyle $ws.Range("A18") $true $true

# -- Row 19 (Skapa funktion...) --
$ws.Range("F19").Value2 = $ws.Range("E19").Value2
$ws.Range("E19").Value2 = $ws.Range("D19").Value2
$ws.Range("D19").Value2 = $ws.Range("C19").Value2
$ws.Range("C19").Value2 = $ws.Range("B19").Value2
$ws.Range("B19").ClearContents() | Out-Null

# -- Row 20 (Skapa en metod...) --
$ws.Range("F20").Value2 = $ws.Range("E20").Value2
$ws.Range("E20").Value2 = $ws.Range("D20").Value2
$ws.Range("D20").Value2 = $ws.Range("C20").Value2
$ws.Range("C20").Value2 = $ws.Range("B20").Value2
$ws.Range("B20").ClearContents() | Out-Null

# -- Row 21 (Skapa metod för introduktionstext) --
$ws.Range("F21").Value2 = $ws.Range("E21").Value2
$ws.Range("E21").Value2 = $ws.Range("D21").Value2
$ws.Range("D21").Value2 = $ws.Range("C21").Value2
$ws.Range("C21").Value2 = $ws.Range("B21").Value2
$ws.Range("B21").ClearContents() | Out-Null

# -- Row 22 (Formatering) -- status flips to "Färdig"
$ws.Range("F22").Value2 = $ws.Range("E22").Value2   # 9
$ws.Range("E22").Value2 = $ws.Range("D22").Value2   # 5
$ws.Range("D22").Value2 = "Färdig"
$ws.Range("C22").Value2 = $ws.Range("B22").Value2   # Formatering
$ws.Range("B22").ClearContents() | Out-Null

# -- Row 23 (Drag/drop-funktion) --
$ws.Range("E23").Value2 = $ws.Range("D23").Value2   # 3 (Uppskattad tid)
$ws.Range("D23").Value2 = $ws.Range("C23").Value2   # Ej påbörjat (Status)
$ws.Range("C23").Value2 = $ws.Range("B23").Value2   # Drag/drop-funktion
$ws.Range("B23").ClearContents() | Out-Null
$ws.Range("F23").ClearContents() | Out-Null

# -- Row 24 (sum row) --
$ws.Range("F24").Value2 = $ws.Range("E24").Value2   # 19
$ws.Range("E24").ClearContents() | Out-Null
Set-CellStyle $ws.Range("E24") $false $true
$ws.Range("D24").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# New marker cells below the table.
# ---------------------------------------------------------------------
$ws.Range("A25").Value2 = $null
Set-CellStyle $ws.Range("A25") $true $false
$ws.Range("A27").Value2 = $null
Set-CellStyle $ws.Range("A27") $true $false
$ws.Range("A29").Value2 = $null
Set-CellStyle $ws.Range("A29") $true $false
$ws.Range("A31").Value2 = $null
Set-CellStyle $ws.Range("A31") $true $false

# ---------------------------------------------------------------------
# Sheet-level cosmetic changes.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 42.85546875
$ws.Range("D1").Select() | Out-Null
